# export-structure-stats.xlsx: add two missing "orientation" rows
# ("Orientation vers CIAS" / "Autre orientation") just above the existing
# "Répartition des orientations" block, fixing a data gap in the export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 103 was an unused blank gap row; inserting a fresh row there pushes
# the old row 104 ("Répartition des orientations" header) and everything
# below it down by one (104->105, ..., 114->115), while row 102 (already
# blank but styled) and the newly inserted row 103 become the two new
# data rows.
$ws.Rows(103).Insert()

# Fill in the two new labels (row 102 already existed/blank, row 103 is
# the freshly inserted one - both inherit the "s=17" list-item style).
$ws.Range("B102").Value = "Orientation vers CIAS"
$ws.Range("B103").Value = "Autre orientation"

# Match the row height used by the rest of this list (16pt).
$ws.Rows(102).RowHeight = 16
$ws.Rows(103).RowHeight = 16

# Restore the active selection to the edited area.
$ws.Range("B101").Select()
